# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker/period table (rows 16-29) is re-sorted so each worker's rows
# are grouped together (periods descending, newest first), and a brand new
# worker (GUSTAVO ADOLFO GARCIA HEREDIA, CC 1143385979) is appended with his
# two overdue periods (2103, 2102) in two new rows. The footer rows (the
# signature block) shift down accordingly, and the summary counters/totals
# at the top of the sheet are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table used to end at row 29 (the row carrying the bottom border of the
# box). Insert two blank rows just above it so it lands on row 31, leaving
# rows 29-30 free for the new worker's two period rows.
$ws.Rows("29:30").Insert()

# The freshly inserted rows inherit a generic default style; clone the
# formatting of a normal data row (row 28) into them so they match the rest
# of the table (borders/fill/number format), leaving the old bottom-border
# row (now row 31) with its own distinct style untouched.
$ws.Range("B28:J28").Copy()
$ws.Range("B29:J30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Full contents of the data table after the edit: CC / doc number / name /
# period / valor mora / salario basico.
$data = @(
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2311", 37333, 1160000),
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2310", 46400, 1160000),
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2309", 46400, 1160000),
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2308", 46400, 1160000),
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2307", 46400, 1160000),
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2306", 46400, 1160000),
  @("CC", "1128056659", "JESICA RODRIGUEZ TEHERAN", "2305", 46400, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2311", 43307, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2310", 46400, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2309", 46400, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2308", 46400, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2307", 46400, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2306", 46400, 1160000),
  @("CC", "45563613",   "LESLY PATRICIA SALCEDO SAMPAYO", "2305", 46400, 1160000),
  @("CC", "1143385979", "GUSTAVO ADOLFO GARCIA HEREDIA", "2103", 36341, 908526),
  @("CC", "1143385979", "GUSTAVO ADOLFO GARCIA HEREDIA", "2102", 36341, 908526)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# Summary block: total overdue value, worker count, period count.
$ws.Range("E11").Value = 710122
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 9
